$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F ("dSF")
$values = @{
    3  = -1
    4  = 6
    5  = 2
    6  = 3
    7  = 4
    8  = 3
    9  = 7
    10 = 1
    11 = 2
    12 = -3
    13 = 3
    14 = 5
    15 = 3
    16 = 2
    17 = 5
    18 = 3
    19 = 4
    20 = -3
    21 = 1
    22 = 1
    23 = 2
    24 = -1
    25 = 1
    26 = 1
    27 = -3
    28 = -3
    29 = 1
    30 = 1
    32 = 2
    33 = -3
    34 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
